$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(7)
$tr = $shape.TextFrame.TextRange
$tr.Text = "Phase1_activities"
